$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet contains a weekly price table for "Cereza" (cherries) at
# "Mercado Mayorista Lo Valledor de Santiago". This edit adds five new
# observation rows to the table (rows that were previously rows 351-360
# get pushed down), and tweaks nothing else.
#
# Plan (matches the target XML row renumbering exactly):
#   1) Insert 4 blank rows at row 351   -> old rows 351-360 become 355-364
#   2) Insert 1 blank row at row 358    -> (old) rows 354-360, now at
#                                          358-364, become 359-365
#   3) Fill in the 5 new rows: 351, 352, 353, 354, 358
# ------------------------------------------------------------------

$ws.Range("A351:A354").EntireRow.Insert()
$ws.Range("A358:A358").EntireRow.Insert()

function Set-Row($r, $A, $B, $C, $D, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T) {
    $ws.Cells.Item($r, 1).Value = $A
    $ws.Cells.Item($r, 2).Value = $B
    $ws.Cells.Item($r, 3).Value = $C
    $ws.Cells.Item($r, 4).Value = $D
    $ws.Cells.Item($r, 5).Value = $E
    $ws.Cells.Item($r, 6).Value = $F
    $ws.Cells.Item($r, 7).Value = $G
    $ws.Cells.Item($r, 8).Value = $H
    $ws.Cells.Item($r, 9).Value = $I
    $ws.Cells.Item($r, 10).Value = $J
    $ws.Cells.Item($r, 11).Value = $K
    $ws.Cells.Item($r, 12).Value = $L
    $ws.Cells.Item($r, 13).Value = $M
    $ws.Cells.Item($r, 14).Value = $N
    $ws.Cells.Item($r, 15).Value = $O
    $ws.Cells.Item($r, 16).Value = $P
    $ws.Cells.Item($r, 17).Value = $Q
    $ws.Cells.Item($r, 18).Value = $R
    $ws.Cells.Item($r, 19).Value = $S
    $ws.Cells.Item($r, 20).Value = $T
}

# Row 351 - Early Burlat, Segunda
Set-Row 351 6 "Mercado Mayorista Lo Valledor de Santiago" "Metropolitana" `
    44516 13 "Fruta" 100103 "Frutos de hueso (carozo)" 100103001 "Cereza" `
    "Early Burlat" "Segunda" 280 11500 11500 11500 `
    "`$/bandeja 5 kilos" "Provincia de Curicó" 2300 5

# Row 352 - Royal Dawn, Primera
Set-Row 352 6 "Mercado Mayorista Lo Valledor de Santiago" "Metropolitana" `
    44516 13 "Fruta" 100103 "Frutos de hueso (carozo)" 100103001 "Cereza" `
    "Royal Dawn" "Primera" 400 28000 28000 28000 `
    "`$/bandeja 10 kilos" "Región de O'Higgins" 2800 10

# Row 353 - Royal Dawn, Primera
Set-Row 353 6 "Mercado Mayorista Lo Valledor de Santiago" "Metropolitana" `
    44516 13 "Fruta" 100103 "Frutos de hueso (carozo)" 100103001 "Cereza" `
    "Royal Dawn" "Primera" 70 45000 45000 45000 `
    "`$/bandeja 10 kilos" "Provincia de Curicó" 3000 15

# Row 354 - Royal Dawn, Segunda
Set-Row 354 6 "Mercado Mayorista Lo Valledor de Santiago" "Metropolitana" `
    44516 13 "Fruta" 100103 "Frutos de hueso (carozo)" 100103001 "Cereza" `
    "Royal Dawn" "Segunda" 30 30000 30000 30000 `
    "`$/bandeja 10 kilos" "Provincia de Curicó" 2000 15

# Row 358 - Royal Dawn, Primor
Set-Row 358 6 "Mercado Mayorista Lo Valledor de Santiago" "Metropolitana" `
    44509 13 "Fruta" 100103 "Frutos de hueso (carozo)" 100103001 "Cereza" `
    "Royal Dawn" "Primor" 160 30000 35000 32500 `
    "`$/bandeja 10 kilos" "Provincia de Curicó" 3250 10

# Make sure the D column (date) keeps the expected date number format style
# that is already used throughout the column (style index copied from the
# row above on insert, so this is mostly a safety net).
$ws.Range("D351:D354").NumberFormat = $ws.Range("D350").NumberFormat
$ws.Range("D358").NumberFormat = $ws.Range("D357").NumberFormat
